$d = $word.ActiveDocument

# --- Move the _GoBack bookmark from the end of the "...in an email." ---
# --- sentence up to right after "Due: Tomorrow at 5pm".              ---
#
# A zero-length Range positioned exactly at a paragraph-end (right
# before the paragraph mark) cannot be used directly as the bookmark
# anchor, so we temporarily insert a placeholder character, anchor the
# bookmark to that (now non-empty) range, then delete the placeholder.
# Word collapses the bookmark to the deletion point, leaving a clean
# zero-length bookmark exactly where we want it.
$rng = $d.Content
$rng.Find.Execute("Due: Tomorrow at 5pm", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$rng.InsertAfter("Z")
$placeholder = $d.Range($rng.Start, $rng.Start + 1)
$d.Bookmarks.Add("_GoBack", $placeholder)
$placeholder2 = $d.Range($rng.Start, $rng.Start + 1)
$placeholder2.Delete()

# Bookmark names are unique, so re-adding "_GoBack" above automatically
# removed it from its old location (just before the final ".") in the
# "Create a public Git repository..." sentence. All that is left to do
# there is merge the now-adjacent, identically-formatted runs back into
# a single run, which a Find/Replace across the join naturally does.
$rng2 = $d.Content
$rng2.Find.Execute("in an email.", $true, $false, $false, $false, $false, $true, 1, $false, "in an email.", 2) | Out-Null
